$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, derived from the source diff.
$changes = @{
    'D2' = '46.561.78'
    'E2' = '  +4.37%  '
    'D3' = '2.472.93'
    'E3' = '  +2.16%  '
    'D4' = '0.999'
    'E4' = '  -0.06%  '
    'D5' = '322.26'
    'E5' = '  +2.01%  '
    'D6' = '105.67'
    'E6' = '  +4.21%  '
    'E7' = '  +1.52%  '
    'D8' = '0.999'
    'E8' = '  -0.07%  '
    'D9' = '0.541'
    'E9' = '  +3.36%  '
    'D10' = '36.24'
    'E10' = '  +2.22%  '
    'D11' = '0.0815'
    'E11' = '  +1.87%  '
    'E12' = '  +0.52%  '
    'D13' = '18.36'
    'E13' = '  -3.00%  '
    'E14' = '  +2.46%  '
    'D15' = '2.861.79'
    'E15' = '  +2.20%  '
    'D16' = '2.483.69'
    'E16' = '  +1.46%  '
    'D17' = '0.846'
    'E17' = '  +1.74%  '
    'D18' = '46.429.85'
    'E18' = '  +4.41%  '
    'D19' = '12.64'
    'E19' = '  +2.54%  '
    'E20' = '  +1.34%  '
    'E21' = '  +2.10%  '
    'D22' = '70.66'
    'E22' = '  +2.78%  '
    'B23' = 'ImmutableX'
    'C23' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D23' = '2.38'
    'E23' = '  +3.95%  '
    'B24' = 'BitcoinCash'
    'C24' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D24' = '248.74'
    'E24' = '  +3.00%  '
    'E25' = '  +2.17%  '
    'D26' = '26.21'
    'E26' = '  +3.96%  '
    'E27' = '  +0.05%  '
    'D28' = '2.20'
    'E28' = '  -3.76%  '
    'E29' = '  +2.91%  '
    'D30' = '34.63'
    'E30' = '  +3.74%  '
    'D31' = '49.68'
    'E31' = '  +2.46%  '
    'E32' = '  +3.46%  '
    'E33' = '  +1.89%  '
    'D34' = '5.32'
    'E34' = '  +3.02%  '
    'E35' = '  -0.01%  '
    'D36' = '0.0768'
    'E37' = '  +2.85%  '
    'E38' = '  +1.32%  '
    'D39' = '2.96'
    'E39' = '  +2.83%  '
    'D40' = '123.37'
    'E40' = '  +2.43%  '
    'E41' = '  +2.27%  '
    'E42' = '  +1.76%  '
    'D43' = '20.89'
    'E43' = '  -0.68%  '
    'D44' = '0.0294'
    'E44' = '  +1.24%  '
    'D45' = '1.986.10'
    'E45' = '  +2.27%  '
    'E46' = '  +1.56%  '
    'E47' = '  -2.48%  '
    'E48' = '  +9.16%  '
    'D49' = '9.08'
    'E49' = '  -3.68%  '
    'D50' = '5.19'
    'E50' = '  +12.11%  '
    'D51' = '79.33'
    'E51' = '  +5.51%  '
}

foreach ($cellRef in $changes.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "46.561.78")
    # are not reinterpreted as numbers/dates by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$cellRef]
    $cell.Style = "Normal"
}
